$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.980814333333332
$ws.Range("H2").Value = 20.942443
$ws.Range("I2").Value = 0.2573350203399358
$ws.Range("J2").Value = 0.2573350203399358
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.61877133333333
$ws.Range("N2").Value = 40.856314
$ws.Range("O2").Value = 0.3264056993691278
$ws.Range("P2").Value = 0.3264056993691277
$ws.Range("Q2").Value = 95.07011412612243
$ws.Range("R2").Value = 855.6310271351019
$ws.Range("S2").Value = 0.08399561728622545
$ws.Range("T2").Value = 0.08399561728622544

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.980814333333332
$ws.Range("H3").Value = 20.942443
$ws.Range("I3").Value = 0.2573350203399358
$ws.Range("J3").Value = 0.2573350203399358
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.92748633333333
$ws.Range("N3").Value = 32.782459
$ws.Range("O3").Value = 0.2619027613928843
$ws.Range("P3").Value = 0.2619027613928842
$ws.Range("Q3").Value = 76.28275322303743
$ws.Range("R3").Value = 686.544779007337
$ws.Range("S3").Value = 0.06739675243012322
$ws.Range("T3").Value = 0.06739675243012321

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.980814333333332
$ws.Range("H4").Value = 20.942443
$ws.Range("I4").Value = 0.2573350203399358
$ws.Range("J4").Value = 0.2573350203399358
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.492675333333334
$ws.Range("N4").Value = 22.478026
$ws.Range("O4").Value = 0.1795794842620271
$ws.Range("P4").Value = 0.1795794842620271
$ws.Range("Q4").Value = 52.30497536194644
$ws.Range("R4").Value = 470.7447782575179
$ws.Range("S4").Value = 0.04621209023520392
$ws.Range("T4").Value = 0.04621209023520391

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.980814333333332
$ws.Range("H5").Value = 20.942443
$ws.Range("I5").Value = 0.2573350203399358
$ws.Range("J5").Value = 0.2573350203399358
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.684515333333332
$ws.Range("N5").Value = 29.053546
$ws.Range("O5").Value = 0.2321120549759609
$ws.Range("P5").Value = 0.2321120549759609
$ws.Range("Q5").Value = 67.60580345031975
$ws.Range("R5").Value = 608.4522310528779
$ws.Range("S5").Value = 0.0597305603883832
$ws.Range("T5").Value = 0.05973056038838319

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.375361666666667
$ws.Range("H6").Value = 19.126085
$ws.Range("I6").Value = 0.2350161092714131
$ws.Range("J6").Value = 0.2350161092714131
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 13.61877133333333
$ws.Range("N6").Value = 40.856314
$ws.Range("O6").Value = 0.3264056993691278
$ws.Range("P6").Value = 0.3264056993691277
$ws.Range("Q6").Value = 86.82459270563223
$ws.Range("R6").Value = 781.4213343506901
$ws.Range("S6").Value = 0.07671059750974696
$ws.Range("T6").Value = 0.07671059750974694

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.375361666666667
$ws.Range("H7").Value = 19.126085
$ws.Range("I7").Value = 0.2350161092714131
$ws.Range("J7").Value = 0.2350161092714131
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.92748633333333
$ws.Range("N7").Value = 32.782459
$ws.Range("O7").Value = 0.2619027613928843
$ws.Range("P7").Value = 0.2619027613928842
$ws.Range("Q7").Value = 69.66667748255723
$ws.Range("R7").Value = 627.000097343015
$ws.Range("S7").Value = 0.06155136798999493
$ws.Range("T7").Value = 0.06155136798999492

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.375361666666667
$ws.Range("H8").Value = 19.126085
$ws.Range("I8").Value = 0.2350161092714131
$ws.Range("J8").Value = 0.2350161092714131
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.492675333333334
$ws.Range("N8").Value = 22.478026
$ws.Range("O8").Value = 0.1795794842620271
$ws.Range("P8").Value = 0.1795794842620271
$ws.Range("Q8").Value = 47.76851510091223
$ws.Range("R8").Value = 429.91663590821
$ws.Range("S8").Value = 0.04220407169622858
$ws.Range("T8").Value = 0.04220407169622856

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.375361666666667
$ws.Range("H9").Value = 19.126085
$ws.Range("I9").Value = 0.2350161092714131
$ws.Range("J9").Value = 0.2350161092714131
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.684515333333332
$ws.Range("N9").Value = 29.053546
$ws.Range("O9").Value = 0.2321120549759609
$ws.Range("P9").Value = 0.2321120549759609
$ws.Range("Q9").Value = 61.74228781637888
$ws.Range("R9").Value = 555.68059034741
$ws.Range("S9").Value = 0.05455007207544269
$ws.Range("T9").Value = 0.05455007207544268

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.973131
$ws.Range("H10").Value = 17.919393
$ws.Range("I10").Value = 0.2201886075150976
$ws.Range("J10").Value = 0.2201886075150976
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.61877133333333
$ws.Range("N10").Value = 40.856314
$ws.Range("O10").Value = 0.3264056993691278
$ws.Range("P10").Value = 0.3264056993691277
$ws.Range("Q10").Value = 81.34670523304467
$ws.Range("R10").Value = 732.120347097402
$ws.Range("S10").Value = 0.07187081642907982
$ws.Range("T10").Value = 0.0718708164290798

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 5.973131
$ws.Range("H11").Value = 17.919393
$ws.Range("I11").Value = 0.2201886075150976
$ws.Range("J11").Value = 0.2201886075150976
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 10.92748633333333
$ws.Range("N11").Value = 32.782459
$ws.Range("O11").Value = 0.2619027613928843
$ws.Range("P11").Value = 0.2619027613928842
$ws.Range("Q11").Value = 65.27130736970966
$ws.Range("R11").Value = 587.4417663273871
$ws.Range("S11").Value = 0.05766800433545805
$ws.Range("T11").Value = 0.05766800433545804

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 5.973131
$ws.Range("H12").Value = 17.919393
$ws.Range("I12").Value = 0.2201886075150976
$ws.Range("J12").Value = 0.2201886075150976
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 7.492675333333334
$ws.Range("N12").Value = 22.478026
$ws.Range("O12").Value = 0.1795794842620271
$ws.Range("P12").Value = 0.1795794842620271
$ws.Range("Q12").Value = 44.75473130646866
$ws.Range("R12").Value = 402.792581758218
$ws.Range("S12").Value = 0.03954135657793514
$ws.Range("T12").Value = 0.03954135657793512

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 5.973131
$ws.Range("H13").Value = 17.919393
$ws.Range("I13").Value = 0.2201886075150976
$ws.Range("J13").Value = 0.2201886075150976
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 9.684515333333332
$ws.Range("N13").Value = 29.053546
$ws.Range("O13").Value = 0.2321120549759609
$ws.Range("P13").Value = 0.2321120549759609
$ws.Range("Q13").Value = 57.84687875750866
$ws.Range("R13").Value = 520.6219088175779
$ws.Range("S13").Value = 0.05110843017262462
$ws.Range("T13").Value = 0.05110843017262462

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 7.798031999999999
$ws.Range("H14").Value = 23.394096
$ws.Range("I14").Value = 0.2874602628735535
$ws.Range("J14").Value = 0.2874602628735535
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 13.61877133333333
$ws.Range("N14").Value = 40.856314
$ws.Range("O14").Value = 0.3264056993691278
$ws.Range("P14").Value = 0.3264056993691277
$ws.Range("Q14").Value = 106.199614658016
$ws.Range("R14").Value = 955.796531922144
$ws.Range("S14").Value = 0.09382866814407556
$ws.Range("T14").Value = 0.09382866814407555

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 7.798031999999999
$ws.Range("H15").Value = 23.394096
$ws.Range("I15").Value = 0.2874602628735535
$ws.Range("J15").Value = 0.2874602628735535
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 10.92748633333333
$ws.Range("N15").Value = 32.782459
$ws.Range("O15").Value = 0.2619027613928843
$ws.Range("P15").Value = 0.2619027613928842
$ws.Range("Q15").Value = 85.212888106896
$ws.Range("R15").Value = 766.915992962064
$ws.Range("S15").Value = 0.07528663663730807
$ws.Range("T15").Value = 0.07528663663730806

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 7.798031999999999
$ws.Range("H16").Value = 23.394096
$ws.Range("I16").Value = 0.2874602628735535
$ws.Range("J16").Value = 0.2874602628735535
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 7.492675333333334
$ws.Range("N16").Value = 22.478026
$ws.Range("O16").Value = 0.1795794842620271
$ws.Range("P16").Value = 0.1795794842620271
$ws.Range("Q16").Value = 58.42812201494399
$ws.Range("R16").Value = 525.8530981344959
$ws.Range("S16").Value = 0.05162196575265948
$ws.Range("T16").Value = 0.05162196575265946

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 7.798031999999999
$ws.Range("H17").Value = 23.394096
$ws.Range("I17").Value = 0.2874602628735535
$ws.Range("J17").Value = 0.2874602628735535
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 9.684515333333332
$ws.Range("N17").Value = 29.053546
$ws.Range("O17").Value = 0.2321120549759609
$ws.Range("P17").Value = 0.2321120549759609
$ws.Range("Q17").Value = 75.52016047382398
$ws.Range("R17").Value = 679.6814442644159
$ws.Range("S17").Value = 0.06672299233951044
$ws.Range("T17").Value = 0.06672299233951043
